$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row height fix-ups (Excel auto-recalculated these after content changed nearby) ---
$ws.Rows(13).RowHeight = 12.5
$ws.Rows(16).RowHeight = 13

# --- Row 14 (no border, "General"-ish styling) ---
$ws.Range("A14").Value = "10월 8일"
$ws.Range("B14").Value = 0.5
$ws.Range("B14").NumberFormat = "h:mm"
$ws.Range("C14").Value = 0.625
$ws.Range("C14").NumberFormat = "h:mm"
$ws.Range("D14").Value = 0
$ws.Range("D14").NumberFormat = "0"
$ws.Range("E14").Value = 180
$ws.Range("E14").NumberFormat = "0_);[Red](0)"
$ws.Range("F14").Value = "GUI 작성"
$ws.Range("A14:F14").Borders.LineStyle = -4142

# --- Row 15 (bordered) ---
$ws.Range("A15").Value = "10월 13일"
$ws.Range("A15").NumberFormat = "m""월"" d""일"""
$ws.Range("B15").Value = 0.83333333333333337
$ws.Range("B15").NumberFormat = "h:mm"
$ws.Range("C15").Value = 0.91666666666666663
$ws.Range("C15").NumberFormat = "h:mm"
$ws.Range("D15").Value = 0
$ws.Range("D15").NumberFormat = "0"
$ws.Range("E15").Value = 120
$ws.Range("E15").NumberFormat = "0_);[Red](0)"
$ws.Range("F15").Value = "Use case specification 작성 회의 "
$ws.Range("A15:F15").Borders.LineStyle = 1

# --- Row 16 (bordered, rich-text date + Dotum activity) ---
$ws.Range("A16").Value = "10월 24일"
$ws.Range("A16").NumberFormat = "m""월"" d""일"""
$ws.Range("A16").Characters(3,1).Font.Name = "돋움"
$ws.Range("A16").Characters(4,3).Font.Name = "Arial"
$ws.Range("A16").Characters(7,1).Font.Name = "돋움"
$ws.Range("B16").Value = 0.75
$ws.Range("B16").NumberFormat = "h:mm"
$ws.Range("C16").Value = 0.875
$ws.Range("C16").NumberFormat = "h:mm"
$ws.Range("D16").Value = 0
$ws.Range("D16").NumberFormat = "0"
$ws.Range("E16").Value = 180
$ws.Range("E16").NumberFormat = "0_);[Red](0)"
$ws.Range("F16").Value = "Use Case Specification 수정"
$ws.Range("F16").Font.Name = "Dotum"
$ws.Range("A16:F16").Borders.LineStyle = 1

# --- Row 17 (no border) ---
$ws.Range("A17").Value = "10월 26일"
$ws.Range("B17").Value = 0.41666666666666669
$ws.Range("B17").NumberFormat = "h:mm"
$ws.Range("C17").Value = 0.5
$ws.Range("C17").NumberFormat = "h:mm"
$ws.Range("D17").Value = 0
$ws.Range("D17").NumberFormat = "0"
$ws.Range("E17").Value = 120
$ws.Range("E17").NumberFormat = "0_);[Red](0)"
$ws.Range("F17").Value = "GUI 수정"

# --- Row 18 (bordered, rich-text date + rich-text activity) ---
$ws.Range("A18").Value = "10월28일"
$ws.Range("A18").NumberFormat = "m""월"" d""일"""
$ws.Range("A18").Characters(3,1).Font.Name = "돋움"
$ws.Range("A18").Characters(4,2).Font.Name = "Arial"
$ws.Range("A18").Characters(6,1).Font.Name = "돋움"
$ws.Range("B18").Value = 0.75
$ws.Range("B18").NumberFormat = "h:mm"
$ws.Range("C18").Value = 0.83333333333333337
$ws.Range("C18").NumberFormat = "h:mm"
$ws.Range("D18").Value = 0
$ws.Range("D18").NumberFormat = "0"
$ws.Range("E18").Value = 120
$ws.Range("E18").NumberFormat = "0_);[Red](0)"
$ws.Range("F18").Value = "SRS회의"
$ws.Range("F18").Characters(4,2).Font.Name = "돋움"
$ws.Range("A18:F18").Borders.LineStyle = 1

# --- Row 19: only the date text was entered so far ---
$ws.Range("A19").Value = "11월 11일"

# --- Cursor ended up on B19 after typing the date and tabbing across ---
$ws.Range("B19").Select()
